$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix fiscal years for Microsoft (rows 2-4 shift from FY2025/2024/2023
# to FY2024/2023/2022) ---

# Row 2: was FY2025 -> now FY2024
$ws.Range("B2").Value = 2024
$ws.Range("C2").Value = 45473
$ws.Range("D2").Value = 245122
$ws.Range("E2").Value = 88136
$ws.Range("F2").Value = 512163
$ws.Range("G2").Value = 243686
$ws.Range("H2").Value = 118548

# Row 3: was FY2024 -> now FY2023
$ws.Range("B3").Value = 2023
$ws.Range("C3").Value = 45107
$ws.Range("D3").Value = 211915
$ws.Range("E3").Value = 72361
$ws.Range("F3").Value = 411976
$ws.Range("G3").Value = 205753
$ws.Range("H3").Value = 87582

# Row 4: was FY2023 -> now FY2022
$ws.Range("B4").Value = 2022
$ws.Range("C4").Value = 44742
$ws.Range("D4").Value = 198270
$ws.Range("E4").Value = 72738
$ws.Range("F4").Value = 364840
$ws.Range("G4").Value = 198298
$ws.Range("H4").Value = 89035

# --- Chatbot prototype work left the sheet selection on I3 ---
$ws.Range("I3").Select() | Out-Null

# --- Page setup was touched (portrait orientation) while exploring
# print/chatbot export options ---
$ws.PageSetup.Orientation = 1
